$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as plain text in the source sheet (e.g. "71.873.78",
# "152.80", "1.00"). Force Text format before writing so COM does not
# auto-convert these numeric-looking strings into actual numbers and drop
# meaningful trailing zeros / the dotted-thousands formatting.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '71.873.78'
$ws.Range('E2').Value = '  +4.81%  '
$ws.Range('D3').Value = '4.036.25'
$ws.Range('E3').Value = '  +4.70%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('E5').Value = '  +3.42%  '
$ws.Range('D6').Value = '152.80'
$ws.Range('E6').Value = '  +8.41%  '
$ws.Range('E7').Value = '  +14.24%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').Value = '0.762'
$ws.Range('E9').Value = '  +7.02%  '
$ws.Range('D10').Value = '0.174'
$ws.Range('E10').Value = '  +3.91%  '
$ws.Range('E11').Value = '  +3.34%  '
$ws.Range('D12').Value = '48.58'
$ws.Range('E12').Value = '  +16.81%  '
$ws.Range('D13').Value = '10.90'
$ws.Range('E13').Value = '  +5.50%  '
$ws.Range('D14').Value = '4.683.10'
$ws.Range('E14').Value = '  +4.66%  '
$ws.Range('D15').Value = '4.032.43'
$ws.Range('E15').Value = '  +4.17%  '
$ws.Range('D16').Value = '14.37'
$ws.Range('E16').Value = '  +2.03%  '
$ws.Range('D17').Value = '20.71'
$ws.Range('E17').Value = '  -3.17%  '
$ws.Range('E18').Value = '  +1.84%  '
$ws.Range('E19').Value = '  -0.10%  '
$ws.Range('D20').Value = '71.828.63'
$ws.Range('E20').Value = '  +4.72%  '
$ws.Range('D21').Value = '436.19'
$ws.Range('E21').Value = '  +5.19%  '
$ws.Range('D22').Value = '99.91'
$ws.Range('E22').Value = '  +15.42%  '
$ws.Range('E23').Value = '  +3.21%  '
$ws.Range('D24').Value = '4.29'
$ws.Range('E24').Value = '  +7.37%  '
$ws.Range('D25').Value = '14.68'
$ws.Range('E25').Value = '  +4.46%  '
$ws.Range('D26').Value = '11.25'
$ws.Range('E26').Value = '  -6.22%  '
$ws.Range('D27').Value = '10.94'
$ws.Range('E27').Value = '  +4.53%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '37.14'
$ws.Range('E28').Value = '  +4.87%  '
$ws.Range('B29').Value = 'LEO'
$ws.Range('C29').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D29').Value = '5.85'
$ws.Range('E29').Value = '  +3.11%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '3.67'
$ws.Range('E30').Value = '  +31.42%  '
$ws.Range('D31').Value = '13.72'
$ws.Range('E31').Value = '  +3.00%  '
$ws.Range('D32').Value = '0.132'
$ws.Range('E32').Value = '  +6.16%  '
$ws.Range('D33').Value = '681.95'
$ws.Range('E33').Value = '  +0.58%  '
$ws.Range('D34').Value = '6.89'
$ws.Range('E34').Value = '  -0.79%  '
$ws.Range('D35').Value = '67.31'
$ws.Range('E35').Value = '  +0.76%  '
$ws.Range('D36').Value = '43.03'
$ws.Range('E36').Value = '  +9.16%  '
$ws.Range('D37').Value = '0.434'
$ws.Range('E37').Value = '  -3.69%  '
$ws.Range('E38').Value = '  +7.31%  '
$ws.Range('D39').Value = '0.0₃0849'
$ws.Range('E39').Value = '  +0.19%  '
$ws.Range('D40').Value = '3.54'
$ws.Range('E40').Value = '  +11.87%  '
$ws.Range('D41').Value = '3.45'
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('D43').Value = '0.0495'
$ws.Range('E43').Value = '  +4.32%  '
$ws.Range('D44').Value = '1.00'
$ws.Range('E45').Value = '  +7.94%  '
$ws.Range('D46').Value = '2.72'
$ws.Range('E46').Value = '  -5.63%  '
$ws.Range('D47').Value = '3.43'
$ws.Range('E47').Value = '  +0.54%  '
$ws.Range('D48').Value = '9.60'
$ws.Range('E48').Value = '  +9.85%  '
$ws.Range('E49').Value = '  +2.11%  '
$ws.Range('D50').Value = '3.38'
$ws.Range('E50').Value = '  +3.19%  '
$ws.Range('D51').Value = '0.000271'
$ws.Range('E51').Value = '  -3.45%  '
